$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New shared strings must be introduced in this exact order so the   ---
# --- sharedStrings table ends up with the same string order as the     ---
# --- target workbook (unused strings get pruned automatically, new     ---
# --- ones get appended in the order they are first written).           ---
$ws.Range("A13").Value = 'all intensity metrics (pair with temp change)'
$ws.Range("A16").Value = '* Change intensity metrics one at a time through four levels (baseline plus 4x4 = 1+16 = 17 runs), then baseline plus temp change, then all intensity metrics, then all intensity metrics plus temp change. '
$ws.Range("F6").Value = 'perc_20th (final)'
$ws.Range("E6").Value = 'perc_20th (raw)'

$ws.Range("A2").Value = 'Precipitation intensity sensitivity analysis'
$ws.Range("A3").Value = 'To be performed (initially) with a baseline scenario: no change in precipitation volume or temperature. '
$ws.Range("B5").Value = 'Emissions scenario'
$ws.Range("A6").Value = 'Precip intensity metric'
$ws.Range("C6").Value = 'wet season precip (%)'
$ws.Range("D6").Value = '3 day max (%)'
$ws.Range("G6").Value = 'perc_80th'
$ws.Range("H6").Value = 'extreme_dry (%)'
$ws.Range("I6").Value = 'extreme wet (%)'
$ws.Range("I6").Font.Bold = $true

$ws.Range("B7").Value = 'RCP 4.5 min'
$ws.Range("C7").Value = -3.76
$ws.Range("D7").Value = -21.25
$ws.Range("E7").Value = -39.06
$ws.Range("F7").Value = 0.12188
$ws.Range("G7").Value = 0.7044
$ws.Range("H7").Value = 'tbd'
$ws.Range("I7").Value = 'tbd'

$ws.Range("B8").Value = 'RCP 4.5 mean'
$ws.Range("C8").Value = 0.84
$ws.Range("D8").Value = -1.7
$ws.Range("E8").Value = -3.317
$ws.Range("F8").Value = 0.193366
$ws.Range("G8").Value = 0.859078
$ws.Range("H8").Value = 'tbd'
$ws.Range("I8").Value = 'tbd'

$ws.Range("B9").Value = 'RCP 4.5 max'
$ws.Range("C9").Value = 4.46
$ws.Range("D9").Value = 10.65
$ws.Range("E9").Value = 32.6
$ws.Range("F9").Value = 0.2652
$ws.Range("G9").ClearContents()
$ws.Range("H9").Value = 'tbd'
$ws.Range("I9").Value = 'tbd'

$ws.Range("B10").Value = 'RCP 8.5 min'
$ws.Range("C10").Value = 2.98
$ws.Range("D10").Value = -14.04
$ws.Range("E10").Value = -54.76
$ws.Range("F10").Value = 0.09048
$ws.Range("G10").ClearContents()
$ws.Range("H10").Value = 'tbd'
$ws.Range("I10").Value = 'tbd'

$ws.Range("B11").Value = 'RCP 8.5 mean'
$ws.Range("C11").Value = 5.56
$ws.Range("D11").Value = 14.35
$ws.Range("E11").Value = -0.42099999999999949
$ws.Range("F11").Value = 0.199158
$ws.Range("G11").Value = 0.90642
$ws.Range("H11").Value = 'tbd'
$ws.Range("I11").Value = 'tbd'

$ws.Range("B12").Value = 'RCP 8.5 max'
$ws.Range("C12").Value = 12.1
$ws.Range("D12").Value = 70.349999999999994
$ws.Range("E12").Value = 40.479999999999997
$ws.Range("F12").Value = 0.28095999999999999
$ws.Range("G12").Value = 1.1408999999999998
$ws.Range("H12").Value = 'tbd'
$ws.Range("I12").Value = 'tbd'

$ws.Range("B14").Value = 'RCP 8.5 mean'
$ws.Range("B14").Font.Bold = $true
$ws.Range("C14").Value = 5.56
$ws.Range("D14").Value = 14.35
$ws.Range("F14").Value = 0.199158
$ws.Range("G14").Value = 0.90642
$ws.Range("H14").Value = 'tbd'
$ws.Range("I14").Value = 'tbd'

# --- Column width / view cosmetics ---
$ws.Columns.Item(3).ColumnWidth = 19.1640625
$ws.Application.ActiveWindow.Zoom = 97
$ws.Range("G19").Select()

Write-Host "done"
